$wb = $excel.ActiveWorkbook

# Re-order the sheets so that "总计" comes first, followed by "2022-Q2"
# (previously "2022-Q2" was first and "总计" was second).
$summary = $wb.Worksheets.Item("总计")
$summary.Move($wb.Worksheets.Item(1))
